$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New command protocol values for costs B3/B4
$ws.Range("B3").Value = 4195.2
$ws.Range("B4").Value = "-"

# Move active selection to B5 (Timelapse setup row)
$ws.Range("B5").Select()
